# Adapt column header formatting to respective input file names (#7)
#
#  - rename header cells ending in "_old" -> "_FV2404"
#  - rename header cells ending in "_new" -> "_FV2410"
#  - wrap the data range in an Excel Table ("Table1") with autofilter
#  - freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1. rename header row (row 1) ------------------------------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $headerText = $cell.Value2
    if ($headerText -ne $null) {
        if ($headerText.EndsWith("_old")) {
            $cell.Value2 = $headerText.Replace("_old", "_FV2404")
        } elseif ($headerText.EndsWith("_new")) {
            $cell.Value2 = $headerText.Replace("_new", "_FV2410")
        }
    }
}

# --- 2. turn the data range into an Excel Table -----------------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# --- 3. freeze the header row -----------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
